# Update database: roll the 5-year reporting window forward by one fiscal
# year (drop 1396/12, add 1401/12) and refresh the figures per the new
# filings/price-read algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "دوره مالی" (fiscal period) headers, shifted one column left,
#     newest period appended in column H ---
$ws.Cells.Item(8,4).Value2 = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8,5).Value2 = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8,6).Value2 = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8,7).Value2 = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8,8).Value2 = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" (publish date) headers, same rolling shift ---
$ws.Cells.Item(9,4).Value2 = "1399-06-20 (9)"
$ws.Cells.Item(9,5).Value2 = "1400-03-30 (9)"
$ws.Cells.Item(9,6).Value2 = "1401-03-25 (8)"
$ws.Cells.Item(9,7).Value2 = "1402-02-28 (7)"
$ws.Cells.Item(9,8).Value2 = "1402-02-28"

# --- Row 11: فروش (Sales) ---
$ws.Cells.Item(11,4).Value2 = 23155
$ws.Cells.Item(11,5).Value2 = 29553
$ws.Cells.Item(11,6).Value2 = 31537
$ws.Cells.Item(11,7).Value2 = 39919
$ws.Cells.Item(11,8).Value2 = 49861

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---
$ws.Cells.Item(12,4).Value2 = -17910
$ws.Cells.Item(12,5).Value2 = -22017
$ws.Cells.Item(12,6).Value2 = -20622
$ws.Cells.Item(12,7).Value2 = -32646
$ws.Cells.Item(12,8).Value2 = -39648

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Cells.Item(13,4).Value2 = 5246
$ws.Cells.Item(13,5).Value2 = 7536
$ws.Cells.Item(13,6).Value2 = 10915
$ws.Cells.Item(13,7).Value2 = 7274
$ws.Cells.Item(13,8).Value2 = 10213

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Cells.Item(14,4).Value2 = -1063
$ws.Cells.Item(14,5).Value2 = -1226
$ws.Cells.Item(14,6).Value2 = -885
$ws.Cells.Item(14,7).Value2 = -1232
$ws.Cells.Item(14,8).Value2 = -2128

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense) ---
$ws.Cells.Item(16,4).Value2 = -529
$ws.Cells.Item(16,5).Value2 = 24
$ws.Cells.Item(16,6).Value2 = 1246
$ws.Cells.Item(16,7).Value2 = 491
$ws.Cells.Item(16,8).Value2 = 469

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Cells.Item(17,4).Value2 = 3654
$ws.Cells.Item(17,5).Value2 = 6333
$ws.Cells.Item(17,6).Value2 = 11276
$ws.Cells.Item(17,7).Value2 = 6533
$ws.Cells.Item(17,8).Value2 = 8554

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Cells.Item(18,4).Value2 = -1301
$ws.Cells.Item(18,5).Value2 = -976
$ws.Cells.Item(18,6).Value2 = -607
$ws.Cells.Item(18,7).Value2 = -1047
$ws.Cells.Item(18,8).Value2 = -1311

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense) ---
$ws.Cells.Item(19,4).Value2 = 182
$ws.Cells.Item(19,5).Value2 = 294
$ws.Cells.Item(19,6).Value2 = 145
$ws.Cells.Item(19,7).Value2 = 242
$ws.Cells.Item(19,8).Value2 = 426

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ---
$ws.Cells.Item(20,4).Value2 = 2535
$ws.Cells.Item(20,5).Value2 = 5651
$ws.Cells.Item(20,6).Value2 = 10814
$ws.Cells.Item(20,7).Value2 = 5728
$ws.Cells.Item(20,8).Value2 = 7669

# --- Row 21: مالیات (Tax) ---
$ws.Cells.Item(21,4).Value2 = -363
$ws.Cells.Item(21,5).Value2 = -876
$ws.Cells.Item(21,6).Value2 = -1272
$ws.Cells.Item(21,7).Value2 = -1041
$ws.Cells.Item(21,8).Value2 = -347

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops) ---
$ws.Cells.Item(22,4).Value2 = 2173
$ws.Cells.Item(22,5).Value2 = 4775
$ws.Cells.Item(22,6).Value2 = 9542
$ws.Cells.Item(22,7).Value2 = 4687
$ws.Cells.Item(22,8).Value2 = 7322

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Cells.Item(24,4).Value2 = 2173
$ws.Cells.Item(24,5).Value2 = 4775
$ws.Cells.Item(24,6).Value2 = 9542
$ws.Cells.Item(24,7).Value2 = 4687
$ws.Cells.Item(24,8).Value2 = 7322

# --- Row 26: سرمایه (Capital) ---
$ws.Cells.Item(26,4).Value2 = 5931
$ws.Cells.Item(26,5).Value2 = 7016
$ws.Cells.Item(26,6).Value2 = 3980
$ws.Cells.Item(26,7).Value2 = 3411
$ws.Cells.Item(26,8).Value2 = 5384
